$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-CellText 2 4 '33.913.15'
Set-CellText 2 5 '  +0.29%  '
Set-CellText 3 4 '1.773.39'
Set-CellText 3 5 '  +0.33%  '
Set-CellText 5 4 '225.36'
Set-CellText 5 5 '  +2.13%  '
Set-CellText 6 4 '0.552'
Set-CellText 6 5 '  +0.63%  '
Set-CellText 7 5 '  -0.04%  '
Set-CellText 8 4 '31.97'
Set-CellText 8 5 '  +3.74%  '
Set-CellText 9 5 '  +1.96%  '
Set-CellText 10 4 '0.0690'
Set-CellText 10 5 '  -2.05%  '
Set-CellText 11 5 '  +1.70%  '
Set-CellText 12 4 '2.031.04'
Set-CellText 12 5 '  +0.43%  '
Set-CellText 13 2 'WrappedEther'
Set-CellText 13 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-CellText 13 4 '1.786.91'
Set-CellText 13 5 '  +0.98%  '
Set-CellText 14 2 'Chainlink'
Set-CellText 14 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText 14 4 '10.94'
Set-CellText 14 5 '  +4.75%  '
Set-CellText 15 4 '33.875.16'
Set-CellText 15 5 '  +0.09%  '
Set-CellText 16 5 '  -0.94%  '
Set-CellText 17 5 '  -0.83%  '
Set-CellText 18 4 '67.14'
Set-CellText 18 5 '  -0.51%  '
Set-CellText 19 4 '239.95'
Set-CellText 19 5 '  -1.11%  '
Set-CellText 20 4 '0.0₃0777'
Set-CellText 20 5 '  +0.88%  '
Set-CellText 21 5 '  +0.14%  '
Set-CellText 22 4 '10.65'
Set-CellText 22 5 '  +1.66%  '
Set-CellText 23 5 '  +1.03%  '
Set-CellText 24 5 '  -2.55%  '
Set-CellText 25 4 '159.46'
Set-CellText 25 5 '  +1.52%  '
Set-CellText 26 4 '16.13'
Set-CellText 27 5 '  +1.47%  '
Set-CellText 28 5 '  +0.91%  '
Set-CellText 29 5 '  +0.11%  '
Set-CellText 30 5 '  +3.00%  '
Set-CellText 31 4 '0.0510'
Set-CellText 31 5 '  -1.89%  '
Set-CellText 32 4 '3.62'
Set-CellText 32 5 '  -1.74%  '
Set-CellText 33 5 '  +1.48%  '
Set-CellText 34 4 '1.79'
Set-CellText 34 5 '  -0.56%  '
Set-CellText 35 4 '1.385.82'
Set-CellText 35 5 '  -0.45%  '
Set-CellText 36 4 '0.652'
Set-CellText 36 5 '  +3.37%  '
Set-CellText 37 5 '  -0.79%  '
Set-CellText 38 4 '0.0186'
Set-CellText 38 5 '  +0.62%  '
Set-CellText 39 5 '  +6.26%  '
Set-CellText 40 5 '  +0.55%  '
Set-CellText 41 4 '0.906'
Set-CellText 41 5 '  -1.59%  '
Set-CellText 42 2 'Aave'
Set-CellText 42 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText 42 4 '77.50'
Set-CellText 42 5 '  -1.17%  '
Set-CellText 43 2 'MXToken'
Set-CellText 43 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText 43 4 '2.66'
Set-CellText 43 5 '  -0.74%  '
Set-CellText 44 2 'InjectiveProtocol'
Set-CellText 44 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-CellText 44 4 '13.32'
Set-CellText 44 5 '  +13.97%  '
Set-CellText 45 2 'BabyDogeCoin'
Set-CellText 45 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-CellText 45 4 '0.0₆0142'
Set-CellText 45 5 '  +21.76%  '
Set-CellText 46 5 '  +4.32%  '
Set-CellText 47 4 '107.99'
Set-CellText 47 5 '  +4.32%  '
Set-CellText 48 4 '0.0496'
Set-CellText 48 5 '  +1.93%  '
Set-CellText 49 4 '5.84'
Set-CellText 49 5 '  -0.35%  '
Set-CellText 50 4 '1.929.24'
Set-CellText 50 5 '  +1.05%  '
Set-CellText 51 5 '  +0.55%  '
